$d = $word.ActiveDocument

# Locate the "Ver no Jupiter Salvar em pdf Salvar em docx" paragraph.
$findRange = $d.Content
$findRange.Find.Execute("Ver no Jupiter Salvar em pdf Salvar em docx", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$verIndex = $findRange.Paragraphs.Item(1).Index

# Locate the site-footer / copyright paragraph that follows it.
$findRange2 = $d.Content
$findRange2.Find.Execute("Powered by Jekyll and Github pages", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$copyIndex = $findRange2.Paragraphs.Item(1).Index

# The blank paragraph immediately before "Ver no Jupiter..." is removed too,
# so the deletion spans from that blank paragraph through the copyright line.
$startPara = $d.Paragraphs.Item($verIndex - 1)
$endPara = $d.Paragraphs.Item($copyIndex)

$delRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$delRange.Delete()
